# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 3 (pushing all existing data rows
# down by one) and populate it with the latest Damasco price report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 3; this shifts rows 3:43 down
# to 4:44 and extends the used range accordingly (A1:T43 -> A1:T44).
$ws.Rows.Item(3).Insert()

$ws.Cells.Item(3, 1).Value  = 10
$ws.Cells.Item(3, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(3, 3).Value  = "La Araucanía"
$ws.Cells.Item(3, 4).Value  = 44545
$ws.Cells.Item(3, 5).Value  = 9
$ws.Cells.Item(3, 6).Value  = "Fruta"
$ws.Cells.Item(3, 7).Value  = 100103
$ws.Cells.Item(3, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(3, 9).Value  = 100103003
$ws.Cells.Item(3, 10).Value = "Damasco"
$ws.Cells.Item(3, 11).Value = "Castle Brite"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 85
$ws.Cells.Item(3, 14).Value = 17000
$ws.Cells.Item(3, 15).Value = 18000
$ws.Cells.Item(3, 16).Value = 17588
$ws.Cells.Item(3, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(3, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(3, 19).Value = 977
$ws.Cells.Item(3, 20).Value = 18
